$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.315.50'
$ws.Range('E2').Value = '  +4.63%  '
$ws.Range('D3').Value = '2.465.66'
$ws.Range('E3').Value = '  +5.92%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = "'567.09"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.78%  '
$ws.Range('D6').Value = "'143.56"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +9.88%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = "'0.590"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.08%  '
$ws.Range('D9').Value = '2.464.62'
$ws.Range('E9').Value = '  +5.91%  '
$ws.Range('E10').Value = '  +4.32%  '
$ws.Range('D11').Value = "'5.73"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.58%  '
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('E13').Value = '  +5.04%  '
$ws.Range('D14').Value = "'26.44"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +12.19%  '
$ws.Range('D15').Value = '2.908.90'
$ws.Range('E15').Value = '  +6.02%  '
$ws.Range('D16').Value = '63.217.15'
$ws.Range('E16').Value = '  +4.53%  '
$ws.Range('E17').Value = '  +6.75%  '
$ws.Range('D18').Value = '2.468.75'
$ws.Range('E18').Value = '  +6.02%  '
$ws.Range('D19').Value = "'11.26"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.78%  '
$ws.Range('D20').Value = "'341.57"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.72%  '
$ws.Range('D21').Value = "'4.30"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.37%  '
$ws.Range('D22').Value = "'6.80"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.98%  '
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').Value = "'65.61"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.14%  '
$ws.Range('E25').Value = '  +1.91%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').Value = "'1.51"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.51%  '
$ws.Range('D28').Value = "'8.13"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.57%  '
$ws.Range('E29').Value = '  +8.43%  '
$ws.Range('D30').Value = '0.0₃0816'
$ws.Range('E30').Value = '  +11.87%  '
$ws.Range('D31').Value = "'6.84"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +13.71%  '
$ws.Range('E32').Value = '  +6.99%  '
$ws.Range('D33').Value = "'176.13"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.52%  '
$ws.Range('E34').Value = '  +10.32%  '
$ws.Range('E35').Value = '  +4.62%  '
$ws.Range('D36').Value = "'18.94"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.95%  '
$ws.Range('D37').Value = "'371.47"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +15.76%  '
$ws.Range('D38').Value = "'4.47"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.98%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('E41').Value = '  +11.79%  '
$ws.Range('D42').Value = "'40.42"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.21%  '
$ws.Range('D43').Value = "'151.66"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +10.21%  '
$ws.Range('D44').Value = "'3.72"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.71%  '
$ws.Range('D45').Value = "'20.63"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.47%  '
$ws.Range('D46').Value = "'0.599"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.88%  '
$ws.Range('D47').Value = "'0.0964"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.32%  '
$ws.Range('D48').Value = "'0.0520"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.60%  '
$ws.Range('D49').Value = '0.0₆0239'
$ws.Range('E49').Value = '  +8.89%  '
$ws.Range('E50').Value = '  +4.55%  '
$ws.Range('D51').Value = "'18.08"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.55%  '
